$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap rows 286 and 287 (REC <-> PPT) ---

$a286 = $ws.Range("A286").Value2
$b286 = $ws.Range("B286").Value2
$c286 = $ws.Range("C286").Value2
$d286 = $ws.Range("D286").Value2
$e286 = $ws.Range("E286").Value2
$f286 = $ws.Range("F286").Value2
$g286 = $ws.Range("G286").Value2

$a287 = $ws.Range("A287").Value2
$b287 = $ws.Range("B287").Value2
$c287 = $ws.Range("C287").Value2
$d287 = $ws.Range("D287").Value2
$e287 = $ws.Range("E287").Value2
$f287 = $ws.Range("F287").Value2
$g287 = $ws.Range("G287").Value2

$ws.Range("A286").Value2 = $a287
$ws.Range("B286").Value2 = $b287
$ws.Range("C286").Value2 = $c287
$ws.Range("D286").Value2 = $d287
$ws.Range("E286").Value2 = $e287
$ws.Range("F286").Value2 = $f287
$ws.Range("G286").Value2 = $g287

$ws.Range("A287").Value2 = $a286
$ws.Range("B287").Value2 = $b286
$ws.Range("C287").Value2 = $c286
$ws.Range("D287").Value2 = $d286
$ws.Range("E287").Value2 = $e286
$ws.Range("F287").Value2 = $f286
$ws.Range("G287").Value2 = $g286

# --- Step 2: insert a new row at 288 for Stuttgart, Germany (STR) ---

$ws.Rows.Item(288).Insert()

$ws.Range("A288").Value2 = "STR"
$ws.Range("B288").Value2 = "Stuttgart, Germany"
$ws.Range("C288").Value2 = 48.783333
$ws.Range("D288").Value2 = 9.183332999999999
$ws.Range("E288").Value2 = "DE"
$ws.Range("F288").Value2 = "Europe"
$ws.Range("G288").Value2 = "Stuttgart"

# match formatting (bordered / bold style) used by the rest of column A
$ws.Range("A288").Borders.LineStyle = 1
